# Insert a new data row at row 922 on the active sheet, shifting the existing
# rows 922:993 down to 923:994 (dimension grows from A1:T993 to A1:T994), then
# populate the newly inserted row with the new "Favorita De Clapp" record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(922).Insert()

$ws.Range("A922").Value = 10
$ws.Range("B922").Value = "Vega Modelo de Temuco"
$ws.Range("C922").Value = "La Araucanía"
$ws.Range("D922").Value = 44931
$ws.Range("E922").Value = 9
$ws.Range("F922").Value = "Fruta"
$ws.Range("G922").Value = 100104
$ws.Range("H922").Value = "Frutos de pepita"
$ws.Range("I922").Value = 100104005
$ws.Range("J922").Value = "Pera"
$ws.Range("K922").Value = "Favorita De Clapp"
$ws.Range("L922").Value = "Primera"
$ws.Range("M922").Value = 55
$ws.Range("N922").Value = 20000
$ws.Range("O922").Value = 20000
$ws.Range("P922").Value = 20000
$ws.Range("Q922").Value = "`$/bandeja 18 kilos granel"
$ws.Range("R922").Value = "Región de O'Higgins"
$ws.Range("S922").Value = 1111
$ws.Range("T922").Value = 18
